$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J6").Value = 3.7
$ws.Range("W6").Value = 1.71
$ws.Range("F7").Value = 1.29
$ws.Range("G7").Value = 1.3
$ws.Range("H7").Value = 12.5
$ws.Range("I7").Value = 16
$ws.Range("L7").Value = 1.27
$ws.Range("N7").Value = 6
$ws.Range("R7").Value = 1.71
$ws.Range("W7").Value = 4.3
$ws.Range("AA7").Value = 610
$ws.Range("AC7").Value = 15
$ws.Range("AF7").Value = 8.6
$ws.Range("AK7").Value = 13.5
$ws.Range("AN7").Value = 3.85
$ws.Range("S8").Value = 3.6
$ws.Range("AB9").Value = 9.4
$ws.Range("L10").Value = 1.01
$ws.Range("Q10").Value = 1.64
$ws.Range("S10").Value = 2.6
$ws.Range("L11").Value = 1.01
$ws.Range("Q12").Value = 1.98
$ws.Range("U12").Value = 2.02
$ws.Range("AC12").Value = 1000
$ws.Range("N13").Value = 1.27
$ws.Range("T13").Value = 1.01
$ws.Range("U13").Value = 1.01
$ws.Range("L14").Value = 1.01
$ws.Range("Q14").Value = 1.51
$ws.Range("R14").Value = 1.68
$ws.Range("AE14").Value = 310
$ws.Range("AJ14").Value = 9.4
$ws.Range("AK14").Value = 17.5
$ws.Range("H15").Value = 4.4
$ws.Range("N15").Value = 3.15
$ws.Range("O15").Value = 1.29
$ws.Range("I16").Value = 3.7
$ws.Range("J16").Value = 3.6
$ws.Range("L16").Value = 1.01
$ws.Range("V16").Value = 1.37
$ws.Range("W16").Value = 1.8
$ws.Range("H17").Value = 8
$ws.Range("L17").Value = 1.01
$ws.Range("AA17").Value = 400
$ws.Range("AC17").Value = 11.5
$ws.Range("AE17").Value = 190
$ws.Range("AH17").Value = 32
$ws.Range("AM17").Value = 220
$ws.Range("F19").Value = 3.45
$ws.Range("G19").Value = 3.5
$ws.Range("J19").Value = 3.65
$ws.Range("W19").Value = 1.4
$ws.Range("AA19").Value = 29
$ws.Range("G20").Value = 1.32
$ws.Range("H20").Value = 11.5
$ws.Range("J20").Value = 6.2
$ws.Range("K20").Value = 7
$ws.Range("S20").Value = 2.1
$ws.Range("W20").Value = 4.1
$ws.Range("AD20").Value = 1000
$ws.Range("AM20").Value = 140
$ws.Range("I21").Value = 6.6
$ws.Range("V21").Value = 1.19
$ws.Range("R22").Value = 1.97
$ws.Range("S22").Value = 1.89
$ws.Range("AB22").Value = 18.5
$ws.Range("P23").Value = 2.2
